# Apply the data changes described in the commit: update the customer /
# product details shown on the return label ("retourLabel") worksheet.
#
# Mapping of old -> new values (shared strings used by the sheet):
#   vincent             -> mustafa             (cells D9 and D11)
#   Orin 15             -> Northview 1          (cell D12)
#   MD-3652             -> 4342                 (cell D13)
#   AKEMI CS-4377-150   -> AKEMI CS-3355-150    (cell D16)
#   4377                -> 3355                 (cell G22)
# (akemi, cell D22, is unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inpaklabel")

$ws.Range("D9").Value = "mustafa"
$ws.Range("D11").Value = "mustafa"
$ws.Range("D12").Value = "Northview 1"
$ws.Range("D16").Value = "AKEMI CS-3355-150"

# D13 and G22 hold digit-only text (order/article numbers stored as text,
# not numbers) and must stay plain text rather than being auto-converted to
# numeric values when assigned directly via .Value. Route the new value
# through a TEXT() formula, then Copy / PasteSpecial(xlPasteValues = -4163)
# over itself: the pasted-in result stays a string (no renewed numeric
# inference) while the cell's original "General" number format/style is
# left untouched.
$ws.Range("D13").Formula = '=TEXT(4342,"0")'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("G22").Formula = '=TEXT(3355,"0")'
$ws.Range("G22").Copy()
$ws.Range("G22").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0
